# Auto-generated PowerShell Excel COM-interop script
# Applies cell value updates to the 'Yojimbo_Profits' workbook sheets
# (market price / leve profit data refreshed by scheduled runner)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 170430.7
$ws.Range("J17").Value = 173696.23
$ws.Range("L17").Value = 521088.6900000001
$ws.Range("N17").Value = -521424.6900000001
$ws.Range("H70").Value = 1688.2941
$ws.Range("I70").Value = 2834
$ws.Range("J70").Value = 1442.7858
$ws.Range("K70").Value = 8502
$ws.Range("L70").Value = 4328.357400000001
$ws.Range("M70").Value = -8232
$ws.Range("N70").Value = -4868.357400000001
$ws.Range("H73").Value = 1688.2941
$ws.Range("I73").Value = 2834
$ws.Range("J73").Value = 1442.7858
$ws.Range("K73").Value = 8502
$ws.Range("L73").Value = 4328.357400000001
$ws.Range("M73").Value = -7566
$ws.Range("N73").Value = -6200.357400000001
$ws.Range("H132").Value = 3050010
$ws.Range("I132").Value = 3572491.8
$ws.Range("J132").Value = 2200
$ws.Range("K132").Value = 10717475.4
$ws.Range("L132").Value = 6600
$ws.Range("M132").Value = -10714945.4
$ws.Range("N132").Value = -11660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1541.5217
$ws.Range("I61").Value = 1541.5217
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1541.5217
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1329.5217
$ws.Range("N61").ClearContents()
$ws.Range("H122").Value = 3821.1667
$ws.Range("I122").Value = 2479.1428
$ws.Range("J122").Value = 5700
$ws.Range("K122").Value = 7437.428400000001
$ws.Range("L122").Value = 17100
$ws.Range("M122").Value = -4987.428400000001
$ws.Range("N122").Value = -22000
$ws.Range("H132").Value = 2344.6667
$ws.Range("I132").Value = 2262.6858
$ws.Range("J132").Value = 3062
$ws.Range("K132").Value = 6788.057400000001
$ws.Range("L132").Value = 9186
$ws.Range("M132").Value = -4258.057400000001
$ws.Range("N132").Value = -14246
$ws.Range("H136").Value = 1541.5217
$ws.Range("I136").Value = 1541.5217
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4624.5651
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2074.5651
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2917.7693
$ws.Range("I105").Value = 3877.5
$ws.Range("J105").Value = 1382.2
$ws.Range("K105").Value = 3877.5
$ws.Range("L105").Value = 1382.2
$ws.Range("M105").Value = -2130.5
$ws.Range("N105").Value = -4876.2
$ws.Range("H134").Value = 1345.4348
$ws.Range("I134").Value = 1320.2273
$ws.Range("J134").Value = 1900
$ws.Range("K134").Value = 3960.6819
$ws.Range("L134").Value = 5700
$ws.Range("M134").Value = -1425.6819
$ws.Range("N134").Value = -10770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19265.754
$ws.Range("I31").Value = 25229.488
$ws.Range("K31").Value = 25229.488
$ws.Range("M31").Value = -24934.488
$ws.Range("H34").Value = 19265.754
$ws.Range("I34").Value = 25229.488
$ws.Range("K34").Value = 25229.488
$ws.Range("M34").Value = -25027.488
$ws.Range("H68").Value = 30294.867
$ws.Range("J68").Value = 30294.867
$ws.Range("L68").Value = 30294.867
$ws.Range("N68").Value = -31792.867
$ws.Range("H71").Value = 30294.867
$ws.Range("J71").Value = 30294.867
$ws.Range("L71").Value = 90884.601
$ws.Range("N71").Value = -98372.601
$ws.Range("H74").Value = 30280.666
$ws.Range("J74").Value = 30280.666
$ws.Range("L74").Value = 30280.666
$ws.Range("N74").Value = -32028.666
$ws.Range("H75").Value = 30260
$ws.Range("J75").Value = 30260
$ws.Range("L75").Value = 30260
$ws.Range("N75").Value = -32256
$ws.Range("H77").Value = 30280.666
$ws.Range("J77").Value = 30280.666
$ws.Range("L77").Value = 90841.99800000001
$ws.Range("N77").Value = -99577.99800000001
$ws.Range("H78").Value = 30260
$ws.Range("J78").Value = 30260
$ws.Range("L78").Value = 90780
$ws.Range("N78").Value = -100764
$ws.Range("H132").Value = 2582.8333
$ws.Range("I132").Value = 1002.4
$ws.Range("J132").Value = 3711.7144
$ws.Range("K132").Value = 3007.2
$ws.Range("L132").Value = 11135.1432
$ws.Range("M132").Value = -477.1999999999998
$ws.Range("N132").Value = -16195.1432
$ws.Range("H134").Value = 2547.4211
$ws.Range("I134").Value = 1638.2
$ws.Range("J134").Value = 5957
$ws.Range("K134").Value = 4914.6
$ws.Range("L134").Value = 17871
$ws.Range("M134").Value = -2379.6
$ws.Range("N134").Value = -22941

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 11346
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 11346
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 34038
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -35410
$ws.Range("H65").Value = 11346
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 11346
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 102114
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -108978

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()
$ws.Range("H102").Value = 2609.1428
$ws.Range("I102").Value = 2387.2
$ws.Range("J102").Value = 3164
$ws.Range("K102").Value = 2387.2
$ws.Range("L102").Value = 3164
$ws.Range("M102").Value = -765.1999999999998
$ws.Range("N102").Value = -6408
$ws.Range("H122").Value = 2784.5
$ws.Range("I122").Value = 2784.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8353.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5903.5
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2291.5217
$ws.Range("I132").Value = 2033.619
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 6100.857
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -3570.857
$ws.Range("N132").Value = -20058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4678.9287
$ws.Range("I132").Value = 4390.222
$ws.Range("J132").Value = 5198.6
$ws.Range("K132").Value = 13170.666
$ws.Range("L132").Value = 15595.8
$ws.Range("M132").Value = -10640.666
$ws.Range("N132").Value = -20655.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 1000000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 1000000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 1000000
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -1000576
$ws.Range("H62").Value = 3899.5
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 5799
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 5799
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -7047
$ws.Range("H65").Value = 3899.5
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 5799
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 28995
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -35235
$ws.Range("H81").Value = 11060
$ws.Range("I81").Value = 1325
$ws.Range("J81").Value = 50000
$ws.Range("K81").Value = 2650
$ws.Range("L81").Value = 100000
$ws.Range("M81").Value = -1589
$ws.Range("N81").Value = -102122
$ws.Range("H84").Value = 11060
$ws.Range("I84").Value = 1325
$ws.Range("J84").Value = 50000
$ws.Range("K84").Value = 13250
$ws.Range("L84").Value = 500000
$ws.Range("M84").Value = -7946
$ws.Range("N84").Value = -510608
$ws.Range("H97").Value = 32000
$ws.Range("J97").Value = 32000
$ws.Range("L97").Value = 32000
$ws.Range("N97").Value = -33982
$ws.Range("H126").Value = 345760.97
$ws.Range("I126").Value = 476945.38
$ws.Range("J126").Value = 1401.875
$ws.Range("K126").Value = 1430836.14
$ws.Range("L126").Value = 4205.625
$ws.Range("M126").Value = -1428366.14
$ws.Range("N126").Value = -9145.625
$ws.Range("H132").Value = 1844
$ws.Range("I132").Value = 1191.8
$ws.Range("K132").Value = 3575.4
$ws.Range("M132").Value = -1045.4
$ws.Range("H136").Value = 780
$ws.Range("I136").Value = 738.75
$ws.Range("J136").Value = 903.75
$ws.Range("K136").Value = 2216.25
$ws.Range("L136").Value = 2711.25
$ws.Range("M136").Value = 333.75
$ws.Range("N136").Value = -7811.25
